$wb = $excel.ActiveWorkbook

# ---- Caso1 ----
$ws = $wb.Worksheets.Item("Caso1")

# Header row (new columns H, I, J)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows 2-19: update G (AutoML MLP recalculated), add H/I/J (Autogluon, H2O, AutoSklearn)
$ws.Range("G2").Value = 1.001238584518433
$ws.Range("H2").Value = 1.003280758857727
$ws.Range("I2").Value = 1.003484047560248
$ws.Range("J2").Value = 1.002901367843151
$ws.Range("G3").Value = 0.9958935976028442
$ws.Range("H3").Value = 1.001055240631104
$ws.Range("I3").Value = 1.001553293987395
$ws.Range("J3").Value = 1.000766569748521
$ws.Range("G4").Value = 0.9968465566635132
$ws.Range("H4").Value = 0.9998693466186523
$ws.Range("I4").Value = 1.000805936497556
$ws.Range("J4").Value = 0.9997910596430302
$ws.Range("G5").Value = 0.9955615401268005
$ws.Range("H5").Value = 0.9996963739395142
$ws.Range("I5").Value = 1.000786427407855
$ws.Range("J5").Value = 1.000250704586506
$ws.Range("G6").Value = 0.9964499473571777
$ws.Range("H6").Value = 0.9992592334747314
$ws.Range("I6").Value = 1.00060666963648
$ws.Range("J6").Value = 0.9997547771781683
$ws.Range("G7").Value = 0.9935828447341919
$ws.Range("H7").Value = 0.9994989633560181
$ws.Range("I7").Value = 1.000615908367428
$ws.Range("J7").Value = 0.999723594635725
$ws.Range("G8").Value = 0.9941142201423645
$ws.Range("H8").Value = 0.9991647005081177
$ws.Range("I8").Value = 1.000623966182114
$ws.Range("J8").Value = 0.9996358714997768
$ws.Range("G9").Value = 0.9961951375007629
$ws.Range("H9").Value = 0.9992102980613708
$ws.Range("I9").Value = 1.000567256549801
$ws.Range("J9").Value = 0.999712185934186
$ws.Range("G10").Value = 0.9944093227386475
$ws.Range("H10").Value = 0.9992414712905884
$ws.Range("I10").Value = 1.000493090101299
$ws.Range("J10").Value = 0.9995849784463644
$ws.Range("G11").Value = 0.9939709901809692
$ws.Range("H11").Value = 0.9995364546775818
$ws.Range("I11").Value = 1.000494003646859
$ws.Range("J11").Value = 0.9994068574160337
$ws.Range("G12").Value = 0.994949996471405
$ws.Range("H12").Value = 0.9992461800575256
$ws.Range("I12").Value = 1.00047822710606
$ws.Range("J12").Value = 0.9996287804096937
$ws.Range("G13").Value = 0.9933212995529175
$ws.Range("H13").Value = 0.999595582485199
$ws.Range("I13").Value = 1.000482487450236
$ws.Range("J13").Value = 0.999556940048933
$ws.Range("G14").Value = 0.9962669610977173
$ws.Range("H14").Value = 0.9993606209754944
$ws.Range("I14").Value = 1.000493802632064
$ws.Range("J14").Value = 0.9996976014226675
$ws.Range("G15").Value = 0.9987419843673706
$ws.Range("H15").Value = 0.9991539120674133
$ws.Range("I15").Value = 1.000522848843885
$ws.Range("J15").Value = 0.9996849428862333
$ws.Range("G16").Value = 0.9956843852996826
$ws.Range("H16").Value = 0.9997172355651855
$ws.Range("I16").Value = 1.000773477298091
$ws.Range("J16").Value = 0.999820988625288
$ws.Range("G17").Value = 0.9983919858932495
$ws.Range("H17").Value = 0.9995560646057129
$ws.Range("I17").Value = 1.000744147704602
$ws.Range("J17").Value = 0.9997209887951612
$ws.Range("G18").Value = 0.9961323142051697
$ws.Range("H18").Value = 0.9992417693138123
$ws.Range("I18").Value = 1.000772291322431
$ws.Range("J18").Value = 0.9996972884982824
$ws.Range("G19").Value = 0.9968447089195251
$ws.Range("H19").Value = 0.9993466138839722
$ws.Range("I19").Value = 1.000760967691581
$ws.Range("J19").Value = 0.9997383747249842

# Match formatting of the new H:J cells to the existing unstyled G column
$ws.Range("G1:G19").Copy()
$ws.Range("H1:J19").PasteSpecial(-4122)

# ---- Caso2 ----
$ws = $wb.Worksheets.Item("Caso2")

# Header row (new columns H, I, J)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows 2-19: update G (AutoML MLP recalculated), add H/I/J (Autogluon, H2O, AutoSklearn)
$ws.Range("G2").Value = 1.003205180168152
$ws.Range("H2").Value = 1.005021095275879
$ws.Range("I2").Value = 1.005107304103829
$ws.Range("J2").Value = 1.004619222134352
$ws.Range("G3").Value = 0.9979731440544128
$ws.Range("H3").Value = 1.002863526344299
$ws.Range("I3").Value = 1.003199179863398
$ws.Range("J3").Value = 1.002466248348355
$ws.Range("G4").Value = 0.9989657998085022
$ws.Range("H4").Value = 1.002026796340942
$ws.Range("I4").Value = 1.002452998277767
$ws.Range("J4").Value = 1.001692758873105
$ws.Range("G5").Value = 0.9976564645767212
$ws.Range("H5").Value = 1.001934170722961
$ws.Range("I5").Value = 1.002434985929292
$ws.Range("J5").Value = 1.002175932750106
$ws.Range("G6").Value = 0.9985307455062866
$ws.Range("H6").Value = 1.001473784446716
$ws.Range("I6").Value = 1.002254480964568
$ws.Range("J6").Value = 1.001541363075376
$ws.Range("G7").Value = 0.9956198334693909
$ws.Range("H7").Value = 1.001758575439453
$ws.Range("I7").Value = 1.002262960582154
$ws.Range("J7").Value = 1.001488450914621
$ws.Range("G8").Value = 0.9961872696876526
$ws.Range("H8").Value = 1.001400351524353
$ws.Range("I8").Value = 1.002270229083444
$ws.Range("J8").Value = 1.00139981135726
$ws.Range("G9").Value = 0.9983062148094177
$ws.Range("H9").Value = 1.001324653625488
$ws.Range("I9").Value = 1.002213923186295
$ws.Range("J9").Value = 1.001411262899637
$ws.Range("G10").Value = 0.9964657425880432
$ws.Range("H10").Value = 1.001449704170227
$ws.Range("I10").Value = 1.002142945662104
$ws.Range("J10").Value = 1.001303384080529
$ws.Range("G11").Value = 0.996023952960968
$ws.Range("H11").Value = 1.001644849777222
$ws.Range("I11").Value = 1.002142871089475
$ws.Range("J11").Value = 1.001223823055625
$ws.Range("G12").Value = 0.9969708323478699
$ws.Range("H12").Value = 1.001353621482849
$ws.Range("I12").Value = 1.002125092049416
$ws.Range("J12").Value = 1.001350037753582
$ws.Range("G13").Value = 0.9953761696815491
$ws.Range("H13").Value = 1.00177526473999
$ws.Range("I13").Value = 1.002128250577641
$ws.Range("J13").Value = 1.001330073922873
$ws.Range("G14").Value = 0.9982562065124512
$ws.Range("H14").Value = 1.001182556152344
$ws.Range("I14").Value = 1.002125583988249
$ws.Range("J14").Value = 1.001325938850641
$ws.Range("G15").Value = 1.00079345703125
$ws.Range("H15").Value = 1.000994920730591
$ws.Range("I15").Value = 1.00215368243638
$ws.Range("J15").Value = 1.001346942037344
$ws.Range("G16").Value = 0.9977428913116455
$ws.Range("H16").Value = 1.001955628395081
$ws.Range("I16").Value = 1.002421410533614
$ws.Range("J16").Value = 1.001634927466512
$ws.Range("G17").Value = 1.00046718120575
$ws.Range("H17").Value = 1.00189197063446
$ws.Range("I17").Value = 1.002393050375884
$ws.Range("J17").Value = 1.00155190192163
$ws.Range("G18").Value = 0.9982463717460632
$ws.Range("H18").Value = 1.001619577407837
$ws.Range("I18").Value = 1.002420021363951
$ws.Range("J18").Value = 1.001590453088284
$ws.Range("G19").Value = 0.9989312887191772
$ws.Range("H19").Value = 1.001808166503906
$ws.Range("I19").Value = 1.002410476437864
$ws.Range("J19").Value = 1.001613674685359

# Match formatting of the new H:J cells to the existing unstyled G column
$ws.Range("G1:G19").Copy()
$ws.Range("H1:J19").PasteSpecial(-4122)

# ---- Caso3 ----
$ws = $wb.Worksheets.Item("Caso3")

# Header row (new columns H, I, J)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("J1").Value = "Valori stimati autosklearn"

# Data rows 2-19: update G (AutoML MLP recalculated), add H/I/J (Autogluon, H2O, AutoSklearn)
$ws.Range("G2").Value = 1.001989603042603
$ws.Range("H2").Value = 1.003972053527832
$ws.Range("I2").Value = 1.00400034281859
$ws.Range("J2").Value = 1.003459963947535
$ws.Range("G3").Value = 0.9967927932739258
$ws.Range("H3").Value = 1.001769304275513
$ws.Range("I3").Value = 1.0021159122846
$ws.Range("J3").Value = 1.001416083425283
$ws.Range("G4").Value = 0.9977542161941528
$ws.Range("H4").Value = 1.000714182853699
$ws.Range("I4").Value = 1.001376899028359
$ws.Range("J4").Value = 1.000558389350772
$ws.Range("G5").Value = 0.9963975548744202
$ws.Range("H5").Value = 1.000703692436218
$ws.Range("I5").Value = 1.001354881648974
$ws.Range("J5").Value = 1.00092407874763
$ws.Range("G6").Value = 0.9972824454307556
$ws.Range("H6").Value = 1.000296354293823
$ws.Range("I6").Value = 1.00117770473943
$ws.Range("J6").Value = 1.000461710616946
$ws.Range("G7").Value = 0.9944061040878296
$ws.Range("H7").Value = 1.000636339187622
$ws.Range("I7").Value = 1.001188940428211
$ws.Range("J7").Value = 1.000495631247759
$ws.Range("G8").Value = 0.9949374198913574
$ws.Range("H8").Value = 1.00042724609375
$ws.Range("I8").Value = 1.001194619667166
$ws.Range("J8").Value = 1.000398155301809
$ws.Range("G9").Value = 0.997018039226532
$ws.Range("H9").Value = 1.000241756439209
$ws.Range("I9").Value = 1.001139489566476
$ws.Range("J9").Value = 1.000479029491544
$ws.Range("G10").Value = 0.9952759742736816
$ws.Range("H10").Value = 1.00033712387085
$ws.Range("I10").Value = 1.001069040137624
$ws.Range("J10").Value = 1.000368108972907
$ws.Range("G11").Value = 0.9947858452796936
$ws.Range("H11").Value = 1.000495076179504
$ws.Range("I11").Value = 1.001067216533095
$ws.Range("J11").Value = 1.000213293358684
$ws.Range("G12").Value = 0.9957616329193115
$ws.Range("H12").Value = 1.000135540962219
$ws.Range("I12").Value = 1.00105237842364
$ws.Range("J12").Value = 1.000386167317629
$ws.Range("G13").Value = 0.9941908717155457
$ws.Range("H13").Value = 1.000491976737976
$ws.Range("I13").Value = 1.001054627890689
$ws.Range("J13").Value = 1.000412167981267
$ws.Range("G14").Value = 0.9969435930252075
$ws.Range("H14").Value = 1.00047755241394
$ws.Range("I14").Value = 1.00107095208346
$ws.Range("J14").Value = 1.000463346019387
$ws.Range("G15").Value = 0.9995073080062866
$ws.Range("H15").Value = 1.000162124633789
$ws.Range("I15").Value = 1.001102643473518
$ws.Range("J15").Value = 1.000451507046819
$ws.Range("G16").Value = 0.9964755773544312
$ws.Range("H16").Value = 1.000767827033997
$ws.Range("I16").Value = 1.001343539543161
$ws.Range("J16").Value = 1.000530855730176
$ws.Range("G17").Value = 0.9992461800575256
$ws.Range("H17").Value = 1.000812172889709
$ws.Range("I17").Value = 1.001312227043152
$ws.Range("J17").Value = 1.000460054725409
$ws.Range("G18").Value = 0.9969591498374939
$ws.Range("H18").Value = 1.000396013259888
$ws.Range("I18").Value = 1.001341229426844
$ws.Range("J18").Value = 1.000424968078732
$ws.Range("G19").Value = 0.9976725578308105
$ws.Range("H19").Value = 1.000474095344543
$ws.Range("I19").Value = 1.001334083376747
$ws.Range("J19").Value = 1.000515632331371

# Match formatting of the new H:J cells to the existing unstyled G column
$ws.Range("G1:G19").Copy()
$ws.Range("H1:J19").PasteSpecial(-4122)
